$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.735.15"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "1.648.35"
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").Value = "1.878.39"
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.644.22"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.532"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.29%  "
$ws.Range("D17").Value = "26.790.44"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").Value = "0.0₃0748"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.63%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("E21").Value = "  +1.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("E29").Value = "  +2.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0523"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("E32").Value = "  +4.40%  "
$ws.Range("E33").Value = "  +2.67%  "
$ws.Range("D34").Value = "1.283.48"
$ws.Range("E34").Value = "  +10.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.30%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0180"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.30%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.31%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.813"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.77%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.517"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.81%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("E41").Value = "  -1.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.807"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.23%  "
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("D44").Value = "1.788.73"
$ws.Range("E44").Value = "  +1.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("E46").Value = "  +4.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.12%  "
$ws.Range("E48").Value = "  -2.33%  "
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0969"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.34%  "
